$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 803, pushing existing rows 803-871 down to 805-873
$ws.Rows(803).Insert()
$ws.Rows(803).Insert()

# Populate new row 803 (copy of row structure with updated pricing data)
$ws.Cells.Item(803, 1).Value = 9
$ws.Cells.Item(803, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(803, 3).Value = "Metropolitana"
$ws.Cells.Item(803, 4).Value = 44783
$ws.Cells.Item(803, 5).Value = 13
$ws.Cells.Item(803, 6).Value = "Fruta"
$ws.Cells.Item(803, 7).Value = 100104
$ws.Cells.Item(803, 8).Value = "Frutos de pepita"
$ws.Cells.Item(803, 9).Value = 100104005
$ws.Cells.Item(803, 10).Value = "Pera"
$ws.Cells.Item(803, 11).Value = "Packham's Triumph"
$ws.Cells.Item(803, 12).Value = "Especial"
$ws.Cells.Item(803, 13).Value = 350
$ws.Cells.Item(803, 14).Value = 12000
$ws.Cells.Item(803, 15).Value = 12000
$ws.Cells.Item(803, 16).Value = 12000
$ws.Cells.Item(803, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(803, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(803, 19).Value = 667
$ws.Cells.Item(803, 20).Value = 18

# Populate new row 804 (copy of row structure with updated pricing data)
$ws.Cells.Item(804, 1).Value = 9
$ws.Cells.Item(804, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(804, 3).Value = "Metropolitana"
$ws.Cells.Item(804, 4).Value = 44783
$ws.Cells.Item(804, 5).Value = 13
$ws.Cells.Item(804, 6).Value = "Fruta"
$ws.Cells.Item(804, 7).Value = 100104
$ws.Cells.Item(804, 8).Value = "Frutos de pepita"
$ws.Cells.Item(804, 9).Value = 100104005
$ws.Cells.Item(804, 10).Value = "Pera"
$ws.Cells.Item(804, 11).Value = "Packham's Triumph"
$ws.Cells.Item(804, 12).Value = "Primera"
$ws.Cells.Item(804, 13).Value = 330
$ws.Cells.Item(804, 14).Value = 10000
$ws.Cells.Item(804, 15).Value = 10000
$ws.Cells.Item(804, 16).Value = 10000
$ws.Cells.Item(804, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(804, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(804, 19).Value = 556
$ws.Cells.Item(804, 20).Value = 18
